# Applies the crypto-price/volume refresh described by the commit diff.
# Cells in column D ("Price") and E ("Volume(1h)") hold text, not numbers
# (e.g. "44.461.69", "  +1.10%  "). Excel auto-coerces plain numeric-looking
# strings assigned to .Value into real numbers (losing formatting / precision),
# so for any new value that parses as a number we briefly force the cell to
# Text format, assign it, then restore General so no stray formatting remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.498.46"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "2.238.60"
$ws.Range("E4").Value = "  +1.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.45"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.41"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.68"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.20"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.834"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").Value = "2.186.33"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.56"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "44.175.30"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").Value = "0.0₃0953"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.94"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.52"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.13"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.47"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.96"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.96"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.04"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0794"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("E34").Value = "  -5.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.110"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +3.99%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.94"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.76"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0300"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "1.792.51"
$ws.Range("E43").Value = "  +3.86%  "
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.65"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +10.37%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "78.80"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -8.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "70.29"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.71"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.90"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.08"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.50"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.68%  "
